$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths (approximate best-fit widths from the source workbook).
# The stored OOXML <col width> is derived from ColumnWidth via the host's
# internal "characters -> serialized width" conversion, so the inputs below
# are chosen to land as close as possible to the authored widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 6
$ws.Columns.Item(2).ColumnWidth = 20.833333333333332
$ws.Columns.Item(3).ColumnWidth = 22.166666666666668
$ws.Columns.Item(4).ColumnWidth = 17.666666666666668
$ws.Columns.Item(5).ColumnWidth = 18.333333333333332
$ws.Columns.Item(6).ColumnWidth = 6
$ws.Columns.Item(7).ColumnWidth = 8.166666666666666

# ---------------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "consumptie huishoudens"
$ws.Range("C1").Value = "investeringen in woningen"
$ws.Range("D1").Value = "bedrijfsinvesteringen"
$ws.Range("E1").Value = "overheidsbestedingen"
$ws.Range("F1").Value = "uitvoer"
$ws.Range("G1").Value = "bbp-groei"

# ---------------------------------------------------------------------------
# Number format for the data block (growth rates, two decimals)
# ---------------------------------------------------------------------------
$ws.Range("B2:G8").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# Data rows
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 2013
$ws.Range("B2").Value = -0.49686056056014999
$ws.Range("C2").Value = -0.37524195707708302
$ws.Range("D2").Value = 0.150499337528466
$ws.Range("E2").Value = -0.12649590593452301
$ws.Range("F2").Value = 0.71792815077837002
$ws.Range("G2").Value = -0.124857872221618

$ws.Range("A3").Value = 2014
$ws.Range("B3").Value = 0.082888168875883905
$ws.Range("C3").Value = 0.202460452469683
$ws.Range("D3").Value = -0.072549583537743906
$ws.Range("E3").Value = 0.042448300865373503
$ws.Range("F3").Value = 1.16816001528366
$ws.Range("G3").Value = 1.42379459921866

$ws.Range("A4").Value = 2015
$ws.Range("B4").Value = 0.30718553055271802
$ws.Range("C4").Value = 0.45544637757114897
$ws.Range("D4").Value = 0.51660946251314199
$ws.Range("E4").Value = -0.067094840138109804
$ws.Range("F4").Value = 0.74729925769011396
$ws.Range("G4").Value = 1.96066341486858

$ws.Range("A5").Value = 2016
$ws.Range("B5").Value = 0.35769030558478099
$ws.Range("C5").Value = 0.60386539481202595
$ws.Range("D5").Value = 0.42815893114750597
$ws.Range("E5").Value = 0.26491741429794902
$ws.Range("F5").Value = 0.53731049316076895
$ws.Range("G5").Value = 2.1898772987727799

$ws.Range("A6").Value = 2017
$ws.Range("B6").Value = 0.57409637317529605
$ws.Range("C6").Value = 0.371589605181181
$ws.Range("D6").Value = 0.25467584592030101
$ws.Range("E6").Value = 0.320894761656915
$ws.Range("F6").Value = 1.34772092959482
$ws.Range("G6").Value = 2.86674805827385

$ws.Range("A7").Value = 2018
$ws.Range("B7").Value = 0.86
$ws.Range("C7").Value = 0.3
$ws.Range("D7").Value = 0.21
$ws.Range("E7").Value = 0.53
$ws.Range("F7").Value = 0.87
$ws.Range("G7").Value = 2.77

$ws.Range("A8").Value = 2019
$ws.Range("B8").Value = 0.7
$ws.Range("C8").Value = 0.23
$ws.Range("D8").Value = 0.03
$ws.Range("E8").Value = 0.79
$ws.Range("F8").Value = 0.83
$ws.Range("G8").Value = 2.57

# ---------------------------------------------------------------------------
# Selection, matching the authored file's saved cursor position
# ---------------------------------------------------------------------------
$ws.Range("B2").Select() | Out-Null

Write-Output "edit applied"
